# Splitting study parameters: insert two new columns (AVG/STD Offset Radial)
# after "STD mAs" and before "kVp", shift the trailing columns right, and
# append a new "Channels" column at the end. Also populate 4 new data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at F (kVp and everything after shifts right two).
$ws.Range("F1:G1").EntireColumn.Insert()

# New header cells for the inserted offset-radial columns (inherit the bold
# /bordered header style from the insert).
$ws.Range("F1").Value = "AVG Offset Radial (mm)"
$ws.Range("G1").Value = "STD Offset Radial (mm)"

# New trailing "Channels" header - copy the header style from "Model" (M1).
$ws.Range("N1").Value = "Channels"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data rows 2-5.
$data = @(
    @(0, "GNL 10 SLICE", 10, 128.4, 119.748235895148, 44.40787540013459, 4.047788885572326, 100, "Br54d", "CT11_Lung", "D:\Ani Nikoghosyan - CT11 vs CT14\Data\Studie_001\CT11_Lung", "SIEMENS", "SOMATOM Force", 96),
    @(1, "GNL 10 SLICE", 10, 130.9, 122.9613353863726, 44.48638906461505, 4.071368516985725, 100, "Br40d", "CT11_Soft", "D:\Ani Nikoghosyan - CT11 vs CT14\Data\Studie_001\CT11_Soft", "SIEMENS", "SOMATOM Force", 96),
    @(2, "GNL 10 SLICE", 10, 38.9, 18.65180956368577, 6.003067895637413, 3.722010148867418, 120, "Bl56f", "CT14_Lung", "D:\Ani Nikoghosyan - CT11 vs CT14\Data\Studie_001\CT14_Lung", "Siemens Healthineers", "NAEOTOM Alpha", 144),
    @(3, "GNL 10 SLICE", 10, 39.1, 18.96549498431296, 6.498821079299242, 3.741648288920148, 120, "Br40f", "CT14_Soft", "D:\Ani Nikoghosyan - CT11 vs CT14\Data\Studie_001\CT14_Soft", "Siemens Healthineers", "NAEOTOM Alpha", 144)
)

# Column A (Index) carries the bold/bordered header-like style too - copy it
# down from A1 across the new rows before filling in values.
$ws.Range("A1").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $ws.Cells.Item($row, 9).Value = $rec[8]
    $ws.Cells.Item($row, 10).Value = $rec[9]
    $ws.Cells.Item($row, 11).Value = $rec[10]
    $ws.Cells.Item($row, 12).Value = $rec[11]
    $ws.Cells.Item($row, 13).Value = $rec[12]
    $ws.Cells.Item($row, 14).Value = $rec[13]
    $row++
}
